$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 671.5897
$ws.Range("I33").Value = 704.44446
$ws.Range("K33").Value = 704.44446
$ws.Range("M33").Value = -475.44446

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2903.4546
$ws.Range("I64").Value = 2886.7
$ws.Range("J64").Value = 2917.4167
$ws.Range("K64").Value = 2886.7
$ws.Range("L64").Value = 2917.4167
$ws.Range("M64").Value = -2638.7
$ws.Range("N64").Value = -3413.4167

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2903.4546
$ws.Range("I67").Value = 2886.7
$ws.Range("J67").Value = 2917.4167
$ws.Range("K67").Value = 2886.7
$ws.Range("L67").Value = 2917.4167
$ws.Range("M67").Value = -2028.7
$ws.Range("N67").Value = -4633.4167

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2910.5
$ws.Range("I69").Value = 2250
$ws.Range("J69").Value = 2983.889
$ws.Range("K69").Value = 6750
$ws.Range("L69").Value = 8951.667000000001
$ws.Range("M69").Value = -5876
$ws.Range("N69").Value = -10699.667

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 2910.5
$ws.Range("I72").Value = 2250
$ws.Range("J72").Value = 2983.889
$ws.Range("K72").Value = 20250
$ws.Range("L72").Value = 26855.001
$ws.Range("M72").Value = -15882
$ws.Range("N72").Value = -35591.001

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 205386.81
$ws.Range("I76").Value = 446910.4
$ws.Range("J76").Value = 4117.1665
$ws.Range("K76").Value = 446910.4
$ws.Range("L76").Value = 4117.1665
$ws.Range("M76").Value = -446595.4
$ws.Range("N76").Value = -4747.1665

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 205386.81
$ws.Range("I79").Value = 446910.4
$ws.Range("J79").Value = 4117.1665
$ws.Range("K79").Value = 446910.4
$ws.Range("L79").Value = 4117.1665
$ws.Range("M79").Value = -445818.4
$ws.Range("N79").Value = -6301.1665

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7382.227
$ws.Range("I132").Value = 9293
$ws.Range("J132").Value = 4622.222
$ws.Range("K132").Value = 27879
$ws.Range("L132").Value = 13866.666
$ws.Range("M132").Value = -25349
$ws.Range("N132").Value = -18926.666

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2651.9385
$ws.Range("I138").Value = 1340.6111
$ws.Range("K138").Value = 4021.8333
$ws.Range("M138").Value = 1118.1667

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1832.1
$ws.Range("I134").Value = 1096
$ws.Range("J134").Value = 3549.6667
$ws.Range("K134").Value = 3288
$ws.Range("L134").Value = 10649.0001
$ws.Range("M134").Value = -753
$ws.Range("N134").Value = -15719.0001

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2100.5293
$ws.Range("I58").Value = 1585
$ws.Range("J58").Value = 2211
$ws.Range("K58").Value = 1585
$ws.Range("L58").Value = 2211
$ws.Range("M58").Value = -1382
$ws.Range("N58").Value = -2617

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 895
$ws.Range("I105").Value = 895
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 895
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 852
$ws.Range("N105").ClearContents()

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2724.0588
$ws.Range("I132").Value = 3113.875
$ws.Range("J132").Value = 2377.5557
$ws.Range("K132").Value = 9341.625
$ws.Range("L132").Value = 7132.6671
$ws.Range("M132").Value = -6811.625
$ws.Range("N132").Value = -12192.6671

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2100.5293
$ws.Range("I136").Value = 1585
$ws.Range("J136").Value = 2211
$ws.Range("K136").Value = 4755
$ws.Range("L136").Value = 6633
$ws.Range("M136").Value = -2205
$ws.Range("N136").Value = -11733

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 781.15625
$ws.Range("I5").Value = 533.9545000000001
$ws.Range("J5").Value = 1325
$ws.Range("K5").Value = 1601.8635
$ws.Range("L5").Value = 3975
$ws.Range("M5").Value = -1489.8635
$ws.Range("N5").Value = -4199

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 781.15625
$ws.Range("I135").Value = 533.9545000000001
$ws.Range("J135").Value = 1325
$ws.Range("K135").Value = 4805.5905
$ws.Range("L135").Value = 11925
$ws.Range("M135").Value = -2270.5905
$ws.Range("N135").Value = -16995

# GSM row 68
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 27500
$ws.Range("I68").Value = 17000
$ws.Range("J68").Value = 38000
$ws.Range("K68").Value = 17000
$ws.Range("L68").Value = 38000
$ws.Range("M68").Value = -16189
$ws.Range("N68").Value = -39622

# GSM row 71
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H71").Value = 27500
$ws.Range("I71").Value = 17000
$ws.Range("J71").Value = 38000
$ws.Range("K71").Value = 51000
$ws.Range("L71").Value = 114000
$ws.Range("M71").Value = -46944
$ws.Range("N71").Value = -122112

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3472.2222
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 4650
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 13950
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -18890

# GSM row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 36554.332
$ws.Range("J139").Value = 36554.332
$ws.Range("L139").Value = 36554.332
$ws.Range("N139").Value = -46834.332

# LTW row 69
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622

# LTW row 72
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14294280
$ws.Range("I132").Value = 29426264
$ws.Range("J132").Value = 2961
$ws.Range("K132").Value = 88278792
$ws.Range("L132").Value = 8883
$ws.Range("M132").Value = -88276262
$ws.Range("N132").Value = -13943

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 34567
$ws.Range("J140").Value = 34567
$ws.Range("L140").Value = 34567
$ws.Range("N140").Value = -44927

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 104194.07
$ws.Range("I62").Value = 4429.905
$ws.Range("J62").Value = 336977.12
$ws.Range("K62").Value = 4429.905
$ws.Range("L62").Value = 336977.12
$ws.Range("M62").Value = -3805.905
$ws.Range("N62").Value = -338225.12

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 104194.07
$ws.Range("I65").Value = 4429.905
$ws.Range("J65").Value = 336977.12
$ws.Range("K65").Value = 22149.525
$ws.Range("L65").Value = 1684885.6
$ws.Range("M65").Value = -19029.525
$ws.Range("N65").Value = -1691125.6
